$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to be treated as plain text,
# matching the original inlineStr cell typing, so Excel does not silently
# coerce numeric-looking strings (e.g. "559.14", "0.370", "3.49") into
# floating point numbers and lose formatting/precision.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.516.45"
$ws.Range("E2").Value = "  +1.24%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.352.44"
$ws.Range("E3").Value = "  -0.82%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "559.14"
$ws.Range("E5").Value = "  +0.58%  "

# Row 6 - Solana
$ws.Range("D6").Value = "175.74"
$ws.Range("E6").Value = "  +2.79%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +1.67%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.344.29"
$ws.Range("E8").Value = "  -0.92%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  -0.08%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +9.25%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +2.94%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "55.20"
$ws.Range("E12").Value = "  -0.76%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +3.94%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +2.10%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.892.06"
$ws.Range("E15").Value = "  -1.25%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "18.25"
$ws.Range("E16").Value = "  +2.61%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.355.91"
$ws.Range("E17").Value = "  -1.74%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -1.00%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "11.81"
$ws.Range("E19").Value = "  +1.57%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "64.442.60"
$ws.Range("E20").Value = "  +1.11%  "

# Row 21 - Polygon
$ws.Range("E21").Value = "  +0.91%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "461.75"
$ws.Range("E22").Value = "  +13.61%  "

# Row 23 - Toncoin
$ws.Range("D23").Value = "4.88"
$ws.Range("E23").Value = "  +11.38%  "

# Row 24 - PancakeSwap
$ws.Range("D24").Value = "4.09"
$ws.Range("E24").Value = "  +0.62%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "86.10"
$ws.Range("E25").Value = "  +4.52%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "13.49"
$ws.Range("E26").Value = "  +2.47%  "

# Row 27 - RenderToken
$ws.Range("E27").Value = "  +1.33%  "

# Row 28 - ImmutableX
$ws.Range("E28").Value = "  +3.61%  "

# Row 29 - Filecoin
$ws.Range("D29").Value = "8.77"
$ws.Range("E29").Value = "  +1.30%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "30.13"
$ws.Range("E30").Value = "  +2.83%  "

# Row 31 - NEARProtocol
$ws.Range("E31").Value = "  +1.11%  "

# Row 32 - Cosmos
$ws.Range("E32").Value = "  +0.91%  "

# Row 33 - Bittensor
$ws.Range("D33").Value = "579.07"
$ws.Range("E33").Value = "  -1.16%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  +1.63%  "

# Row 35 - OKB
$ws.Range("D35").Value = "59.10"
$ws.Range("E35").Value = "  +1.18%  "

# Row 37 - Kaspa
$ws.Range("E37").Value = "  -6.11%  "

# Row 38 & 39 - swap InjectiveProtocol/Stacks rows, with updated values
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "3.49"
$ws.Range("E38").Value = "  +3.28%  "

$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").Value = "35.82"
$ws.Range("E39").Value = "  +0.09%  "

# Row 40 - PEPE
$ws.Range("D40").Value = "0.0₃0753"
$ws.Range("E40").Value = "  +4.27%  "

# Row 41 - TheGraph
$ws.Range("D41").Value = "0.370"
$ws.Range("E41").Value = "  +0.80%  "

# Row 42 - Maker
$ws.Range("D42").Value = "3.085.13"
$ws.Range("E42").Value = "  -2.72%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  -0.24%  "

# Row 44 - ThetaToken
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  -0.60%  "

# Row 45 - Fetch.AI
$ws.Range("D45").Value = "2.50"
$ws.Range("E45").Value = "  +0.48%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  +2.77%  "

# Row 47 - ApeXProtocol
$ws.Range("D47").Value = "3.21"
$ws.Range("E47").Value = "  +0.28%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  +2.68%  "

# Row 49 - WEMIXToken
$ws.Range("D49").Value = "2.59"
$ws.Range("E49").Value = "  -0.32%  "

# Row 50 & 51 - swap THORChain/Monero rows, with updated values
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "136.34"
$ws.Range("E50").Value = "  +1.33%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "8.34"
$ws.Range("E51").Value = "  +2.28%  "
